$wb = $excel.ActiveWorkbook
$wsInputs = $wb.Worksheets.Item("Inputs")
$wsOutputs = $wb.Worksheets.Item("Outputs")

# --- Inputs sheet edits ---

# D2 comment simplified
$wsInputs.Range("D2").Value = "class 1"

# B3 / B4 formulas simplified to their literal computed values
$wsInputs.Range("B3").Formula = "=8839"
$wsInputs.Range("B4").Formula = "=52.778"

# --- Outputs sheet edits ---

# Fix typo "diagra" -> "diagram" on existing rows
$wsOutputs.Range("D2").Value = "Wing loading diagram"
$wsOutputs.Range("D3").Value = "Wing loading diagram"
$wsOutputs.Range("D4").Value = "Wing loading diagram"

# Widen column D
$wsOutputs.Columns.Item(4).ColumnWidth = 23.17

# New row 5: Wing loading at cruise
$wsOutputs.Range("A5").Value = "Wing loading at cruise"
$wsOutputs.Range("B5").Value = 4405.636769211125
$wsOutputs.Range("C5").Value = "N/m^2"
$wsOutputs.Range("D5").Value = "Wing loading diagram"

# New row 6: Wing loading at landing
$wsOutputs.Range("A6").Value = "Wing loading at landing"
$wsOutputs.Range("B6").Value = 4872.47953908665
$wsOutputs.Range("C6").Value = "N/m^2"
$wsOutputs.Range("D6").Value = "Wing loading diagram"

# --- Selection / active sheet updates ---
# Outputs keeps F4 as its stored selection (not the active tab anymore)
$wsOutputs.Activate() | Out-Null
$wsOutputs.Range("F4").Select() | Out-Null

# Inputs becomes the active (selected) tab, with B3 selected - activate last
# so it ends up as the workbook's active sheet.
$wsInputs.Activate() | Out-Null
$wsInputs.Range("B3").Select() | Out-Null
